# #1088 changed all strings in input .xlsx files to lower case and underscores
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val.ToString().ToLower().Replace(" ", "_")
        $cell.Value = $newVal
    }
}

# Move selection to A15, mirroring the post-edit cursor position recorded in the file
$ws.Range("A15").Select()
